# Applies the "Updated cryptos list" price/volume refresh described by the diff.
# Most D-column price strings look like plain text (e.g. thousand-dot formatted
# "65.610.08"), but some of the refreshed values are lone decimals (e.g. "571.80")
# that Excel would otherwise auto-convert to a Number. A leading apostrophe forces
# those into Text cells, matching the original workbook formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.610.08'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.952.63'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''571.80'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = '''162.35'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '''0.517'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '2.949.10'
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('D10').Value = '''6.70'
$ws.Range('E10').Value = '  -3.76%  '
$ws.Range('E11').Value = '  -4.14%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '''0.0000244'
$ws.Range('E13').Value = '  -2.72%  '
$ws.Range('D14').Value = '''34.82'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '65.619.90'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '3.440.39'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''15.98'
$ws.Range('E19').Value = '  +14.28%  '
$ws.Range('B20').Value = 'WrappedEther'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D20').Value = '2.950.34'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').Value = '''445.75'
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('D22').Value = '''0.696'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = '''82.10'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('E25').Value = '  -3.36%  '
$ws.Range('D26').Value = '''12.30'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('D27').Value = '''10.07'
$ws.Range('E27').Value = '  -4.80%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').Value = '''2.53'
$ws.Range('E29').Value = '  +8.14%  '
$ws.Range('D30').Value = '''8.07'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''2.60'
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '''0.0000102'
$ws.Range('E32').Value = '  -0.77%  '
$ws.Range('E33').Value = '  +3.26%  '
$ws.Range('D34').Value = '''27.23'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '''0.973'
$ws.Range('E36').Value = '  -2.15%  '
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('D38').Value = '''45.76'
$ws.Range('E38').Value = '  +4.76%  '
$ws.Range('D39').Value = '''49.20'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  -7.25%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.122'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = '''0.302'
$ws.Range('E42').Value = '  -2.62%  '
$ws.Range('D43').Value = '''2.84'
$ws.Range('E43').Value = '  -6.36%  '
$ws.Range('D44').Value = '''8.55'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = '''385.55'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('E46').Value = '  -1.59%  '
$ws.Range('D47').Value = '2.679.40'
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('D48').Value = '''133.20'
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('D50').Value = '''23.82'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('E51').Value = '  +0.77%  '
